# Apply the IG-regeneration style update:
#  - bump the "pythia" IG name to "cicada" in the embedded URLs
#  - bump the generation Date
#  - add a new "Jurisdiction" metadata row (blank value) before "Description"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata"
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL row (row 2): pythia -> cicada
$meta.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/preferred-vaccine-reason"

# Date row (row 8): new generation timestamp
$meta.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" row before the existing "Description" row (row 11),
# copying the formatting of the preceding row so the inserted row matches style.
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$meta.Cells.Item(11, 1).Value = "Jurisdiction"

# ---------------------------------------------------------------------------
# Sheet "Elements"
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.url fixed value (row 5, column R = "Fixed Value"): pythia -> cicada
$elements.Range("R5").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/preferred-vaccine-reason"

# Extension.value[x] binding value set (row 6, column Z = "Binding Value Set"): pythia -> cicada
$elements.Range("Z6").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/preferred-allowed-reason"
